# "Generate Report for Archive"
# - Update the localization status from "Ready for handoff" to "In Translation"
#   on every sheet that reports it (Overview!E2:F2, zh-cn!C2, de-de!C2).
# - Narrow the now-shorter status columns to match (Overview!E:F, zh-cn!C, de-de!C).

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"

# Overview sheet: zh-cn / de-de status cells live in columns E and F (row 2)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2:F2").Value = $newStatus
$wsOverview.Range("E1:F1").ColumnWidth = 12.5

# zh-cn sheet: Status column is C (row 2)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C1").ColumnWidth = 12.5

# de-de sheet: Status column is C (row 2)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C1").ColumnWidth = 12.5
